# Update the "as_of_utc" timestamp column (AA) on the data sheets.
# The diff shows every data row (2-26) on both the "Главные" and
# "Линейные" sheets getting their AA column value updated from
# "2025-11-26 03:06:02" to "2025-11-26 07:06:43". These are stored as
# plain text strings (inline strings), not Excel date values, so we
# must assign them as text to avoid Excel converting them into a
# numeric date serial.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-11-26 03:06:02"
$newTimestamp = "2025-11-26 07:06:43"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    if ($lastRow -lt 2) { $lastRow = 2 }

    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA = 27
        if ($cell.Value2 -eq $oldTimestamp) {
            $cell.Value2 = $newTimestamp
        }
    }
}
